$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.296.24'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.930.53'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7466'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '249.88'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9992'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3222'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '28.00'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07125'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7892'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08011'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.934.58'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.394'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.50'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.50'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.302.35'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '253.12'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008070'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.768'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.184.91'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.836'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.588'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.12'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.09'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1335'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.295'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.359'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.531'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.421'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.296'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.05114'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7478'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.769'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01979'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.800'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.36'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.400'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4498'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.991'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8441'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9992'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.48'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.799'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.528'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '999.30'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +13.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.30'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06064'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.28%  '
